$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column G holds header "K" (Strike#) values that need to be regenerated.
$ws.Range("G2").Value = 5
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 5
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 3
